$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix swapped "Support Contents" values between the burstReload row (9)
# --- and the switchableFiremodes row (6): they had been copy/pasted into the
# --- wrong rows. Also add the missing "Default" value for switchableFiremodes
# --- and a "Notes" remark for burstReload.

# Row 6 (switchableFiremodes): Support Contents should be the firemode list,
# and Default should be "Auto".
$ws.Range("F6").Value = "Auto, Semi, Burst, Safe"
$ws.Range("G6").Value = "Auto"

# F6 previously carried a (now incorrect) top border inherited from the
# header-like formatting; clear it so the cell matches its neighbours.
$ws.Range("F6").Borders.Item(3).LineStyle = -4142

# Row 9 (burstReload): Support Contents goes back to "/" and a clarifying
# note is added.
$ws.Range("F9").Value = "/"
$ws.Range("H9").Value = "Only works in Burst firemode"

# --- Update the saved view state (scroll position/selection) ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E15").Select()
